$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "MuSCs"
$ws.Range("M2").Value = 0.01393633333333333
$ws.Range("N2").Value = 0.041809
$ws.Range("Q2").Value = 0.007770490667555556
$ws.Range("R2").Value = 0.069934416008
